$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Cases" Cypher query shown in B2 dropped its trailing
# "Cohort" column (coalesce(co.cohort_description, '') AS `Cohort`)
# and lost a stray blank line after the first MATCH clause.
$ws.Range("B2").Value = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC02'] and diag.stage_of_disease in ['T3N0M0', 'T3N1M0', 'T3N1M1'] OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# Shorter query text means the wrapped rows no longer need to be quite
# so tall.
$ws.Rows.Item(2).RowHeight = 290
$ws.Rows.Item(3).RowHeight = 290
$ws.Rows.Item(4).RowHeight = 290

# Re-point the view/selection at the top of the sheet, on B2.
$ws.Activate()
$ws.Range("B2").Select()
